$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 171, shifting existing rows 171-216 down to 172-217
$ws.Rows.Item(171).Insert()

# Populate the newly inserted row 171 with the new data record
$ws.Cells.Item(171, 1).Value = 4
$ws.Cells.Item(171, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(171, 3).Value = "Los Lagos"
$ws.Cells.Item(171, 4).Value = 44508
$ws.Cells.Item(171, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(171, 5).Value = 10
$ws.Cells.Item(171, 6).Value = 100112008
$ws.Cells.Item(171, 7).Value = "Coliflor"
$ws.Cells.Item(171, 8).Value = "Sin especificar"
$ws.Cells.Item(171, 9).Value = "Primera"
$ws.Cells.Item(171, 10).Value = 500
$ws.Cells.Item(171, 11).Value = 1200
$ws.Cells.Item(171, 12).Value = 1200
$ws.Cells.Item(171, 13).Value = 1200
$ws.Cells.Item(171, 14).Value = "$/unidad"
$ws.Cells.Item(171, 15).Value = "Región Metropolitana"
$ws.Cells.Item(171, 16).Value = 1200
$ws.Cells.Item(171, 17).Value = 1
$ws.Cells.Item(171, 18).Value = "Hortaliza"
